# Insert a new data row at row 77 (pushing the existing rows 77-148 down to
# 78-149) and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 77; all rows below shift down by one.
$ws.Rows.Item(77).Insert()

# Make the new date cell (D77) use the same number format as the date column
# elsewhere (e.g. the cell that just moved into D78) before setting its value.
$ws.Range("D77").NumberFormat = $ws.Range("D78").NumberFormat

$ws.Range("A77").Value = 9
$ws.Range("B77").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C77").Value = "Metropolitana"
$ws.Range("D77").Value = 44902
$ws.Range("E77").Value = 13
$ws.Range("F77").Value = 100112022
$ws.Range("G77").Value = "Arveja Verde"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 32
$ws.Range("K77").Value = 22000
$ws.Range("L77").Value = 23000
$ws.Range("M77").Value = 22375
$ws.Range("N77").Value = "$/saco 25 kilos"
$ws.Range("O77").Value = "Provincia de Talca"
$ws.Range("P77").Value = 895
$ws.Range("Q77").Value = 25
$ws.Range("R77").Value = "Hortaliza"
